$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.815.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  +3.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.59%  "

$ws.Range("E6").Value = "  +2.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4429"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.90%  "

$ws.Range("E9").Value = "  +2.93%  "

$ws.Range("E10").Value = "  +1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.885.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.565"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.765"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07233"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.21%  "

$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009163"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.036"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  +2.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.844.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.332"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.971"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.992"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.331"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09108"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7807"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.219"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.585"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.52%  "

$ws.Range("E35").Value = "  +3.27%  "

$ws.Range("E36").Value = "  +2.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05363"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.856"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5211"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.892"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.703"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.729"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.65%  "

$ws.Range("E47").Value = "  +3.06%  "

$ws.Range("E48").Value = "  +3.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.905"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.55%  "
